$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B3").Value = "3 : résultats postés ou publiés après les 36 mois"
$ws.Range("D3").Value = "2013-002869-19"
$ws.Range("B4").Value = "2 : résultats postés ou publiés entre 12 et 36 mois"
$ws.Range("B5").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "2012-004555-36"
$ws.Range("G5").Value = "A prospective randomized controlled multicentre trial comparing half-dose photodynamic therapy (PDT) with high-density subthreshold micropulse laser treatment in patients with chronic central serous chorioretinopathy (CSC)."
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "DRUG (presumed)"
$ws.Range("B6").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("A7").Value = "4"
$ws.Range("B7").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C7").Value = "NCT02899806"
$ws.Range("F7").Value = "2017"
$ws.Range("G7").Value = "Impact of a Video Explaining Epidural Analgesia in Obstetrics in Terms of Satisfaction, Understanding and Anxiety: A Prospective Randomised Trial"
$ws.Range("H7").Value = "VIDEOCLIP"
$ws.Range("I7").Value = "OTHER"
$ws.Range("B8").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B9").Value = "1 : résultats postés ou publiés dans les 12 mois"
$ws.Range("B10").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C10").Value = "NCT01490580"
$ws.Range("F10").Value = "2020"
$ws.Range("G10").Value = "Double Blind Randomized Controlled Trial Comparing ""Atropine+Propofol"" Versus ""Atropine+Atracurium+Sufentanil"" as a Premedication Prior to Semi-urgent or Elective Endotracheal Intubation of Term and Preterm Newborns"
$ws.Range("H10").Value = "PRETTINEO"
$ws.Range("I10").Value = "DRUG"
$ws.Range("B11").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C11").Value = "NCT04776174"
$ws.Range("G11").Value = "Crossover Comparison of the Efficacy and Tolerance of Telerobotic vs Standard Ultrasound Exam in Children"
$ws.Range("H11").Value = ""
$ws.Range("A12").Value = "4"
$ws.Range("B12").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C12").Value = "NCT03939377"
$ws.Range("G12").Value = "Evaluation of an Osteopathic Procedure in the Management of Pain in Palliative Care Patients in a Mobile Palliative Care Support Team (EMASP): Controlled, Randomized, Single-blind Study"
$ws.Range("H12").Value = "OSTEOPAL"
$ws.Range("I12").Value = "OTHER"
$ws.Range("B13").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C13").Value = "NCT03803228"
$ws.Range("G13").Value = "Comparison of the Cumulative Number of Oocytes Obtained With 2 Controlled Ovarian Hyperstimulations (COH) Within the Same Cycle With FertistartKit® (DUOSTIM) Versus 2 Conventional COH in Poor Ovarian Responders Undergoing IVF. Bistim Study"
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = "DRUG"
$ws.Range("B14").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C14").Value = "NCT04068558"
$ws.Range("G14").Value = "Synchronized Nasal Intermittent Positive Pressure Ventilation Versus Noninvasive Neurally Adjusted Ventilatory Assist Ventilation in Extremely Premature Infants: a Randomized Crossover Trial"
$ws.Range("H14").Value = "EASYNNEO"
$ws.Range("I14").Value = "DEVICE"
$ws.Range("B15").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C15").Value = "NCT03488758"
$ws.Range("G15").Value = "Enjoyment of Infant Formulas Based on Cow or Goat Milk Protein"
$ws.Range("H15").Value = "CHARLIE"
$ws.Range("I15").Value = "OTHER"
$ws.Range("A16").Value = "2"
$ws.Range("B16").Value = "2 : résultats postés ou publiés entre 12 et 36 mois"
$ws.Range("C16").Value = "NCT03030664"
$ws.Range("F16").Value = "2021"
$ws.Range("G16").Value = "Randomised Controlled Trial With Two Parallel Arms Testing the Effect of L. Reuteri on Bowel Movements in Children Aged 6 Months to 4 Years"
$ws.Range("H16").Value = "BIOWELL"
$ws.Range("I16").Value = "DIETARY_SUPPLEMENT"
$ws.Range("B17").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C17").Value = "NCT05079139"
$ws.Range("F17").Value = "2022"
$ws.Range("G17").Value = "Musset's Surgical Technique: Evaluation of Long-term Results (LONGOMUSSET)"
$ws.Range("H17").Value = "LONGOMUSSET"
$ws.Range("I17").Value = "PROCEDURE"
$ws.Range("B18").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B19").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C19").Value = "NCT03895099"
$ws.Range("G19").Value = "Feasibility and Efficacy of a New Ovarian Stimulation Regimen With RANDom Start, Use of Corifollitropin Alpha and Progestin Protocol for Oocyte donorS"
$ws.Range("H19").Value = "RANDOS"
$ws.Range("B20").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("A21").Value = "4"
$ws.Range("B21").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C21").Value = "NCT02884245"
$ws.Range("F21").Value = "2023"
$ws.Range("G21").Value = "Interest of Estrogen Scheduling Before Ovarian Stimulation With Corifollitropin Alfa in Women Older Than 38 Years Old Undergoing in Vitro Fertilization"
$ws.Range("H21").Value = "PRESCORI"
$ws.Range("I21").Value = "DRUG"
$ws.Range("A22").Value = "4"
$ws.Range("B22").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C22").Value = "NCT05092659"
$ws.Range("F22").Value = "2025"
$ws.Range("G22").Value = "Patient Experience in Bariatric Surgery: Definition of New Indicators"
$ws.Range("H22").Value = "CALEX"
$ws.Range("I22").Value = "BEHAVIORAL"
$ws.Range("A23").Value = "4"
$ws.Range("B23").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("D23").Value = "2009-011403-23"
$ws.Range("G23").Value = "ETUDE PILOTE SUR L’EFFICACITE ET LA TOLERANCE DES INJECTIONS INTRA-VITREENNES DE LUCENTIS (RANIBIZUMABÒ) A LA PHASE INITIALE DES OCCLUSIONS DE LA VEINE CENTRALE DE LA RETINE"
$ws.Range("I23").Value = "DRUG (presumed)"
